$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Station Delft" row (row 26) into the small table (columns A-D)
# (Shared-string insertion order matters: B26 must be written before A26)
$ws.Range("B26").Value = "Station Delft, 2611 AC Delft"
$ws.Range("A26").Value = "Station Delft"
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = "begin_eindpunt"

# Match formatting of the rest of the table (A16:A25 style)
$ws.Range("A26").Style = $ws.Range("A25").Style

# Adjust column D width (per diff: bestFit removed, custom width 13.90625)
$ws.Columns("D").ColumnWidth = 13.90625

# Update view/selection state to match target workbook
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("E33").Select()
